$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the polite_expressions value on the existing row 7 (it moves down to row 8)
$ws.Cells.Item(7, 3).Value = ""

# Insert a new row 8 with the new review entry, shifting the old row 7's
# "polite_expressions" value ("nan") down while updating the rest of the
# columns for the new review.
$ws.Cells.Item(8, 1).Value = "parisk"
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = "nan"
$ws.Cells.Item(8, 4).Value = "DIS"
$ws.Cells.Item(8, 5).Value = "EXP"
$ws.Cells.Item(8, 6).Value = "d3fb2dcb-ee08-4432-9f4b-c252dbb3433f"
$ws.Cells.Item(8, 7).Value = "SJ3dBGZ0Z_annotated.xlsx"
$ws.Cells.Item(8, 8).Value = "We evaluate our method on NLP task for two reasons: 1) they are particularly well-suited for evaluating our method (naturally large output spaces) 2) we did not dispose of the computational resources to tackle tasks from other domains such as vision (e.g. Flickr100M) which requires hundreds of GPUs for weeks."
$ws.Cells.Item(8, 9).Value = "Correct"
